$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out all existing data on the sheet (rows 1-6 plus any stray formatting
# on rows 7/18) so only the single remaining business exception survives.
$ws.Cells.Clear()

# Leave only the one business-exception row that should remain.
$ws.Range("A1").Value = "The Conclusion Evidence Location path: \\MainFolder\Remediation_or_Justification Evidence\1-22-2020\CHR0000391114.pdf for the change: CHANGES - SOX Audit Report for magic_qq_23455.txt_07.01.73.eml made on 1/22/2020 is not a valid path."
